# Refresh the cryptocurrency price table (rows 2-51) with the latest
# scraped values from coinranking.com, matching the GitHub Actions run.
#
# Price cells (column D) are stored as plain text in this sheet (many values
# such as "26.207.33" or "1.658.35" are not valid numbers anyway, but some
# refreshed prices *are* valid-looking numbers, e.g. "0.5178"). A leading
# apostrophe forces Excel to keep those as text instead of silently
# converting them to floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.207.33"
$ws.Range("E2").Value = "  -4.05%  "

# Row 3
$ws.Range("D3").Value = "1.658.35"
$ws.Range("E3").Value = "  -2.82%  "

# Row 4
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$ws.Range("E5").Value = "  -2.54%  "

# Row 6
$ws.Range("D6").Value = "'0.5178"
$ws.Range("E6").Value = "  -2.77%  "

# Row 7
$ws.Range("D7").Value = "'1.006"
$ws.Range("E7").Value = "  +0.34%  "

# Row 8
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "'0.06429"
$ws.Range("E8").Value = "  -2.63%  "

# Row 9
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.2565"
$ws.Range("E9").Value = "  -4.08%  "

# Row 10
$ws.Range("D10").Value = "'19.94"
$ws.Range("E10").Value = "  -4.77%  "

# Row 11
$ws.Range("D11").Value = "'0.07789"

# Row 12
$ws.Range("D12").Value = "1.655.06"
$ws.Range("E12").Value = "  -3.24%  "

# Row 13
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.886.55"
$ws.Range("E13").Value = "  -2.94%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.296"
$ws.Range("E14").Value = "  -5.54%  "

# Row 15
$ws.Range("D15").Value = "'0.5534"
$ws.Range("E15").Value = "  -3.99%  "

# Row 16
$ws.Range("D16").Value = "0.0₅8062"
$ws.Range("E16").Value = "  -1.34%  "

# Row 17
$ws.Range("D17").Value = "'64.43"
$ws.Range("E17").Value = "  -4.81%  "

# Row 18
$ws.Range("D18").Value = "26.222.51"
$ws.Range("E18").Value = "  -3.95%  "

# Row 19
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'212.01"
$ws.Range("E19").Value = "  -2.35%  "

# Row 20
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.006"
$ws.Range("E20").Value = "  +0.36%  "

# Row 21
$ws.Range("D21").Value = "'4.388"
$ws.Range("E21").Value = "  -5.91%  "

# Row 22
$ws.Range("E22").Value = "  -3.41%  "

# Row 23
$ws.Range("D23").Value = "'5.909"
$ws.Range("E23").Value = "  -0.89%  "

# Row 24
$ws.Range("D24").Value = "'1.006"
$ws.Range("E24").Value = "  +0.31%  "

# Row 25
$ws.Range("D25").Value = "'143.42"
$ws.Range("E25").Value = "  +1.06%  "

# Row 26
$ws.Range("D26").Value = "'1.757"
$ws.Range("E26").Value = "  +1.16%  "

# Row 27
$ws.Range("D27").Value = "'0.1165"
$ws.Range("E27").Value = "  -3.92%  "

# Row 28
$ws.Range("D28").Value = "'6.973"
$ws.Range("E28").Value = "  -3.85%  "

# Row 29
$ws.Range("D29").Value = "'15.78"
$ws.Range("E29").Value = "  -2.77%  "

# Row 30
$ws.Range("D30").Value = "'0.05282"
$ws.Range("E30").Value = "  -1.99%  "

# Row 31
$ws.Range("E31").Value = "  -2.52%  "

# Row 32
$ws.Range("D32").Value = "'3.369"
$ws.Range("E32").Value = "  -3.68%  "

# Row 33
$ws.Range("D33").Value = "'3.225"
$ws.Range("E33").Value = "  -5.80%  "

# Row 34
$ws.Range("D34").Value = "'1.575"
$ws.Range("E34").Value = "  -4.14%  "

# Row 35
$ws.Range("D35").Value = "'2.763"
$ws.Range("E35").Value = "  -3.83%  "

# Row 36
$ws.Range("D36").Value = "'2.363"
$ws.Range("E36").Value = "  -2.02%  "

# Row 37
$ws.Range("E37").Value = "  -2.18%  "

# Row 38
$ws.Range("D38").Value = "'0.5716"
$ws.Range("E38").Value = "  -2.31%  "

# Row 39
$ws.Range("D39").Value = "1.167.35"
$ws.Range("E39").Value = "  +11.63%  "

# Row 40
$ws.Range("E40").Value = "  -2.58%  "

# Row 41
$ws.Range("D41").Value = "'1.006"
$ws.Range("E41").Value = "  +0.35%  "

# Row 42
$ws.Range("D42").Value = "'0.8380"
$ws.Range("E42").Value = "  -0.24%  "

# Row 43
$ws.Range("D43").Value = "'5.669"
$ws.Range("E43").Value = "  -3.06%  "

# Row 44
$ws.Range("D44").Value = "'99.96"
$ws.Range("E44").Value = "  -0.91%  "

# Row 45
$ws.Range("D45").Value = "1.796.26"
$ws.Range("E45").Value = "  -3.01%  "

# Row 46
$ws.Range("D46").Value = "0.0₈109"
$ws.Range("E46").Value = "  -1.03%  "

# Row 47
$ws.Range("D47").Value = "'0.4503"
$ws.Range("E47").Value = "  -0.22%  "

# Row 48
$ws.Range("D48").Value = "'56.04"
$ws.Range("E48").Value = "  -3.31%  "

# Row 49
$ws.Range("D49").Value = "'1.010"
$ws.Range("E49").Value = "  +0.37%  "

# Row 50
$ws.Range("D50").Value = "'7.900"
$ws.Range("E50").Value = "  -2.24%  "

# Row 51
$ws.Range("D51").Value = "'0.05084"
$ws.Range("E51").Value = "  -2.81%  "
